# sample.xlsx changes: reorder "lecturer_attendance" sheet to sit right after
# "schedule" (before "topic"), add a new schedule_id column + a 3rd data row to
# that sheet, update the schedule sheet's start/end time values, and move the
# active/selected tab from "fic" to "schedule".

$wb = $excel.ActiveWorkbook

# --- 1. Move "lecturer_attendance" tab to just after "schedule" -------------
$src    = $wb.Worksheets.Item("lecturer_attendance")
$target = $wb.Worksheets.Item("schedule")
$src.Move($null, $target) | Out-Null

# --- 2. Extend "lecturer_attendance" with a schedule_id column + new row ----
$la = $wb.Worksheets.Item("lecturer_attendance")

$la.Range("G1").Value = "schedule_id"
$la.Range("G2").Value = 1
$la.Range("G3").Value = 2

$la.Range("A4").Value = 3
$la.Range("B4").Value = 1518166800
$la.Range("C4").Value = 1518166800
$la.Range("D4").Value = 1518174000
$la.Range("E4").Value = 1
$la.Range("F4").Value = 3
$la.Range("G4").Value = 3

# --- 3. Update "schedule" sheet start/end time values -----------------------
$sched = $wb.Worksheets.Item("schedule")

$sched.Range("B2").Value = 1517958000
$sched.Range("C2").Value = 1517968200
$sched.Range("B3").Value = 1518066000
$sched.Range("C3").Value = 1518073200
$sched.Range("B4").Value = 1518166800
$sched.Range("C4").Value = 1518174000
$sched.Range("B5").Value = 1518231600
$sched.Range("C5").Value = 1518238800

# --- 4. Move the active-tab / selection from "fic" to "schedule" -----------
$la.Range("E9").Select() | Out-Null
$sched.Activate() | Out-Null
$sched.Range("F9").Select() | Out-Null
